# Swap the presentation's active theme colour scheme (the slide master's
# theme, which is the one PowerPoint actually renders/exports) so it takes
# on the "Office Theme" colours that used to live in the deck's other,
# unused theme part. (dk1/lt1 are black/white in both themes already, but
# they're set here too for completeness.)
#
# ThemeColorScheme.Colors(index).RGB uses the classic OLE colour encoding
# (0x00BBGGRR), i.e. bytes in the opposite order from a plain #RRGGBB
# value, so colours are built from R/G/B components rather than literal
# hex triples to avoid byte-order mistakes.
function Set-ThemeColor($colorObj, [int]$r, [int]$g, [int]$b) {
    $colorObj.RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
Set-ThemeColor $scheme.Colors(1)  0x00 0x00 0x00   # dk1
Set-ThemeColor $scheme.Colors(2)  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $scheme.Colors(3)  0x44 0x54 0x6A   # dk2
Set-ThemeColor $scheme.Colors(4)  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $scheme.Colors(5)  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $scheme.Colors(6)  0xED 0x7D 0x31   # accent2
Set-ThemeColor $scheme.Colors(7)  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $scheme.Colors(8)  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $scheme.Colors(9)  0x44 0x72 0xC4   # accent5
Set-ThemeColor $scheme.Colors(10) 0x70 0xAD 0x47   # accent6
Set-ThemeColor $scheme.Colors(11) 0x05 0x63 0xC1   # hlink
Set-ThemeColor $scheme.Colors(12) 0x95 0x4F 0x72   # folHlink
